# Game of Thrones answer key workbook update.
# Fills in the final-season answers (columns D/E/F) for questions 23-48
# (rows 25-50) and corrects the "Lives"/"Dies" outcome for several
# characters in rows 2-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Flip "Lives" -> "Dies" for the characters who do not survive ---
$diesRows = @(3, 9, 11, 12, 14, 16, 19, 20, 23)
foreach ($r in $diesRows) {
    $ws.Cells.Item($r, 6).Value = "Dies"
}

# --- New answers for rows 25-49 (also flip the "correctAnswer" flag D to TRUE) ---
# Cell values are set in the same order the strings were first introduced so the
# shared-string table is built up identically to the authored workbook.
$ws.Cells.Item(3, 6).Value = "Dies"
$ws.Cells.Item(32, 6).Value = "No"
$ws.Cells.Item(31, 6).Value = "Edd Tollett"
$ws.Cells.Item(45, 6).Value = "Rhaegal"
$ws.Cells.Item(28, 6).Value = "The Hound kills Ser Gregor, Ser Gregor kills Qyburn"
$ws.Cells.Item(25, 6).Value = "No One/There will be no Iron Throne"
$ws.Cells.Item(29, 6).Value = "There is no prince that was promised"
$ws.Cells.Item(42, 6).Value = "Edmure Tully, Robin Arryn"
$ws.Cells.Item(44, 6).Value = "Nymeria, Ghost, Drogon"
$ws.Cells.Item(48, 6).Value = "Gilly and Samwell "
$ws.Cells.Item(27, 6).Value = "Arya kills the Night's King"

# --- Remaining F-column answers that reuse already-existing shared strings ---
$ws.Cells.Item(26, 6).Value = "Sansa Stark"
$ws.Cells.Item(30, 6).Value = "Arya Stark"
$ws.Cells.Item(33, 6).Value = "No"
$ws.Cells.Item(34, 6).Value = "No"
$ws.Cells.Item(35, 6).Value = "No"
$ws.Cells.Item(36, 6).Value = "No"
$ws.Cells.Item(37, 6).Value = "No"
$ws.Cells.Item(38, 6).Value = "No"
$ws.Cells.Item(39, 6).Value = "No"
$ws.Cells.Item(41, 6).Value = "No"
$ws.Cells.Item(46, 6).Value = "No"
$ws.Cells.Item(47, 6).Value = "No"
$ws.Cells.Item(49, 6).Value = "Jon Snow"

# Row 50's answer is numeric (number of sips), not a shared string.
$ws.Cells.Item(50, 6).Value = 15

# --- Apply the "answer" style (same look as F2/F9/...) to the cells that
#     carry it in the target file ---
$styledRows = @(25, 26, 27, 29)
foreach ($r in $styledRows) {
    $ws.Cells.Item($r, 6).Font.Color = 0
}

# --- Mark every one of these newly-answered questions as "correct" (column D) ---
$correctRows = @(25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 41, 42, 44, 45, 46, 47, 48, 49)
foreach ($r in $correctRows) {
    $ws.Cells.Item($r, 4).Value = $true
}

# --- Update the sort-state range to cover the new F column, and restore the
#     previously-active selection. ---
$ws.Range("F28").Select()

"done"
